# working_hours.xlsx edit:
# Insert a new blank separator row before the summary rows, and fill the
# previously-empty "D144" row with an actual working-hours entry
# (2014-07-21, start time 0.76388888888888884).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the blank row (old row 144) and the three summary rows
# (old rows 145-147) down by one, so a new blank row appears at 145
# and the summary rows become 146-148. Excel automatically fixes up the
# SUM/F145/F146 formulas to reference the new row numbers.
$ws.Rows("145:145").Insert()

# Populate the now-freed row 144 with the new data entry.
$ws.Range("A144").Value = 2014
$ws.Range("B144").Value = 7
$ws.Range("C144").Value = 21
$ws.Range("D144").Value = 0.76388888888888884

# Match the author's final selection (cell F144).
$ws.Range("F144").Select()
